$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "SBO_DEF" header in column F, matching the formatting of the
# existing header cells (bold font + border, like B1:E1)
$ws.Range("F1").Value = "SBO_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Fill the new column's data rows (F2:F10) with the placeholder value "[]"
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
